# 2017/9/19 meeting minutes update
#
# 1) Slide 5 speaker notes: the "2017/9/19 ... undo operation" line gets
#    "future work：" inserted before "刪除" (the "刪除後需要可以救回，" run is
#    also split into "刪除" + "後需要可以救回，", same text, two runs).
# 2) Remove all the leftover "弧形箭號 (上彎/下彎)" (curved arrow) decoration
#    shapes from every slide (slides 1-6).

$p = $ppt.ActivePresentation

# --- 1) Update slide 5's notes text -----------------------------------
$notesSlide = $p.Slides.Item(5).NotesPage
for ($i = 1; $i -le $notesSlide.Shapes.Count; $i++) {
    $shp = $notesSlide.Shapes.Item($i)
    if ($shp.HasTextFrame) {
        $tr = $shp.TextFrame.TextRange
        if ($tr.Text -like "*2017/9/19*" -and $tr.Text -like "*undo*") {
            $tr.Text = "2017/9/19 future work：刪除後需要可以救回，undo operation"
        }
    }
}

# --- 2) Delete the curved-arrow decoration shapes on every slide -----
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    $toDelete = @()
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shp = $slide.Shapes.Item($j)
        if ($shp.Name -like "*弧形箭號*") {
            $toDelete += $shp
        }
    }
    foreach ($shp in $toDelete) {
        $shp.Delete()
    }
}
